$wb = $excel.ActiveWorkbook

# Insert a new worksheet "TestSheet_1_withMW" right after "TestSheet_1" that keeps
# a full copy of the original data (including the MW column), then strip the MW
# column values from the original "TestSheet_1" sheet (keeping the cell styling).

$orig = $wb.Worksheets.Item("TestSheet_1")
$copy = $wb.Worksheets.Add($null, $orig)
$copy.Name = "TestSheet_1_withMW"

$orig.Range("A1:N13").Copy($copy.Range("A1"))

$copy.Range("A1").Select()

$orig.Range("I2:I13").ClearContents()
